$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 duplicates row 2's data (same venue/date/result/teams/player/stats).
# Columns G:K hold digit-only strings ("2", "4", "0", "0", "50.00") that must
# stay text (matching row 2's t="str" cells) instead of being auto-converted
# to numbers - format them as Text before assigning so "50.00" keeps its
# trailing zeros and the column stays consistent with the existing data.
$ws.Range("G3:K3").NumberFormat = "@"

$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " September 19 2020"
$ws.Range("C3").Value = "Super Kings won by 5 wickets (with 4 balls remaining)"
$ws.Range("D3").Value = "Mumbai Indians"
$ws.Range("E3").Value = "Chennai Super Kings"
$ws.Range("F3").Value = "Rahul Chahar" + [char]0x00A0   # trailing char matches F2's non-breaking space
$ws.Range("G3").Value = "2"
$ws.Range("H3").Value = "4"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "50.00"
